$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ordering of (id, speaker_variant) pairs for rows 2-10,
# with is_prefered (column D) cleared for every row.
$rows = @(
    @{ B = "#attalia";      C = "Attalia" },
    @{ B = "#sardanapalus"; C = "Sardanapalus" },
    @{ B = "#oporus";       C = "Oporus" },
    @{ B = "#hypermis";     C = "Hypermis" },
    @{ B = "#belus";        C = "Belus" },
    @{ B = "#arbaces";      C = "Arbaces" },
    @{ B = "#porus";        C = "Porus" },
    @{ B = "#salomena";     C = "Salomena" },
    @{ B = "#ninias";       C = "Ninias" }
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $null
    $r++
}
